$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the sample order data rows (2-6) with the generic field-name
# placeholders now used by OrderManager's database-backed lookup.
$headers = @("sku", "name", "quantity", "cost_per", "total_cost")

for ($row = 2; $row -le 6; $row++) {
    for ($col = 1; $col -le 5; $col++) {
        $ws.Cells.Item($row, $col).Value = $headers[$col - 1]
    }
}
